$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.336.79"
$ws.Range("E2").Value = "  +2.02%  "
$ws.Range("D3").Value = "2.975.80"
$ws.Range("E3").Value = "  +2.03%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'503.80"
$ws.Range("E5").Value = "  +7.95%  "
$ws.Range("D6").Value = "'134.59"
$ws.Range("E6").Value = "  +9.09%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +7.37%  "
$ws.Range("D9").Value = "'7.31"
$ws.Range("E9").Value = "  +12.59%  "
$ws.Range("E10").Value = "  +12.64%  "
$ws.Range("E11").Value = "  +8.63%  "
$ws.Range("E12").Value = "  +3.64%  "
$ws.Range("D13").Value = "3.486.69"
$ws.Range("E13").Value = "  +2.01%  "
$ws.Range("D14").Value = "'25.15"
$ws.Range("E14").Value = "  +12.64%  "
$ws.Range("D15").Value = "'0.0000152"
$ws.Range("E15").Value = "  +15.75%  "
$ws.Range("D16").Value = "56.361.27"
$ws.Range("E16").Value = "  +2.03%  "
$ws.Range("D17").Value = "2.976.28"
$ws.Range("E17").Value = "  +1.83%  "
$ws.Range("D18").Value = "'5.69"
$ws.Range("E18").Value = "  +13.31%  "
$ws.Range("D19").Value = "'12.33"
$ws.Range("E19").Value = "  +9.09%  "
$ws.Range("D20").Value = "'7.75"
$ws.Range("E20").Value = "  +11.61%  "
$ws.Range("D21").Value = "'324.08"
$ws.Range("E21").Value = "  +6.53%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "'0.470"
$ws.Range("E23").Value = "  +7.22%  "
$ws.Range("D24").Value = "'61.91"
$ws.Range("E24").Value = "  +5.64%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").Value = "'0.162"
$ws.Range("E26").Value = "  +7.34%  "
$ws.Range("D27").Value = "0.0₃0885"
$ws.Range("E27").Value = "  +12.30%  "
$ws.Range("D28").Value = "'6.48"
$ws.Range("E28").Value = "  +10.62%  "
$ws.Range("D29").Value = "'6.78"
$ws.Range("E29").Value = "  +14.11%  "
$ws.Range("E30").Value = "  +8.09%  "
$ws.Range("E31").Value = "  +12.21%  "
$ws.Range("D32").Value = "'20.49"
$ws.Range("E32").Value = "  +8.77%  "
$ws.Range("D33").Value = "'158.33"
$ws.Range("E33").Value = "  +12.11%  "
$ws.Range("E34").Value = "  +8.55%  "
$ws.Range("D35").Value = "'1.26"
$ws.Range("E35").Value = "  +6.86%  "
$ws.Range("D36").Value = "'5.55"
$ws.Range("E36").Value = "  +4.86%  "
$ws.Range("E37").Value = "  +12.33%  "
$ws.Range("D38").Value = "'22.89"
$ws.Range("E38").Value = "  +8.14%  "
$ws.Range("D39").Value = "3.009.84"
$ws.Range("E39").Value = "  +2.10%  "
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("D41").Value = "'36.22"
$ws.Range("E41").Value = "  +4.12%  "
$ws.Range("D42").Value = "'0.639"
$ws.Range("E42").Value = "  +6.97%  "
$ws.Range("D43").Value = "2.245.15"
$ws.Range("E43").Value = "  +10.91%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "'1.39"
$ws.Range("E44").Value = "  +9.23%  "
$ws.Range("B45").Value = "ONDO"
$ws.Range("C45").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D45").Value = "'0.982"
$ws.Range("E45").Value = "  +4.22%  "
$ws.Range("E46").Value = "  +6.36%  "
$ws.Range("D47").Value = "'1.93"
$ws.Range("E47").Value = "  +27.63%  "
$ws.Range("E48").Value = "  +14.08%  "
$ws.Range("E49").Value = "  +11.64%  "
$ws.Range("D50").Value = "'18.97"
$ws.Range("E50").Value = "  +9.89%  "
$ws.Range("D51").Value = "'0.0864"
$ws.Range("E51").Value = "  +10.38%  "
